# report_3.xlsx — refresh the sample "Dépenses" row with new fixture data
# (004: fix mvn job)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text fields — safe to assign directly
$ws.Range("A2").Value = "shopping with Escher"
$ws.Range("E2").Value = "Durable Aluminum Pants"
$ws.Range("F2").Value = "Pandanus leaves"

# Fields whose new text looks like a date / number to Excel's parser —
# force the Text number format first so the literal string is kept
# instead of being auto-converted to a date serial / numeric value.
$dateCell = $ws.Range("B2")
$dateCell.NumberFormat = "@"
$dateCell.Value = "2003-02-05"

$priceCell = $ws.Range("G2")
$priceCell.NumberFormat = "@"
$priceCell.Value = "4.0"

$qtyCell = $ws.Range("H2")
$qtyCell.NumberFormat = "@"
$qtyCell.Value = "4.0"
